$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings used only by A8 / C9) ---
$ws.Range("A8").Value = "Volume 30   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# --- Numeric cell updates ---
$ws.Range("L14").Value = -83.333333333333
$ws.Range("I15").Value = 11
$ws.Range("K15").Value = 10
$ws.Range("L15").Value = 175
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -70.27027027027
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -71.428571428571
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -52.631578947368
$ws.Range("I16").Value = 77
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = 24.193548387096
$ws.Range("L16").Value = 24.193548387096
$ws.Range("M16").Value = -26.666666666666
$ws.Range("N16").Value = -83.297180043383
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 80
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 116
$ws.Range("J17").Value = 119
$ws.Range("K17").Value = -2.521008403361
$ws.Range("L17").Value = -3.333333333333
$ws.Range("M17").Value = 6.422018348623
$ws.Range("N17").Value = -62.58064516129
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -23.529411764705
$ws.Range("I18").Value = 69
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = -23.333333333333
$ws.Range("L18").Value = 16.949152542372
$ws.Range("M18").Value = -2.81690140845
$ws.Range("N18").Value = -74.444444444444
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -35.135135135135
$ws.Range("I19").Value = 131
$ws.Range("J19").Value = 147
$ws.Range("K19").Value = -10.884353741496
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 29.702970297029
$ws.Range("N19").Value = 5.645161290322
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -12.5
$ws.Range("J20").Value = 59
$ws.Range("K20").Value = -35.593220338983
$ws.Range("L20").Value = 35.714285714285
$ws.Range("M20").Value = 2.702702702702
$ws.Range("N20").Value = -84.166666666666
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -15.384615384615
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -24.107142857142
$ws.Range("I21").Value = 443
$ws.Range("J21").Value = 490
$ws.Range("K21").Value = -9.591836734693
$ws.Range("L21").Value = 8.048780487804
$ws.Range("M21").Value = 0.226244343891
$ws.Range("N21").Value = -69.490358126721
$ws.Range("L22").Value = -12.5
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = -28.571428571428
$ws.Range("I23").Value = 32
$ws.Range("J23").Value = 27
$ws.Range("K23").Value = 18.518518518518
$ws.Range("L23").Value = -23.809523809523
$ws.Range("M23").Value = 3.225806451612
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 77.777777777777
$ws.Range("F24").Value = 66
$ws.Range("G24").Value = 49
$ws.Range("H24").Value = 34.69387755102
$ws.Range("I24").Value = 320
$ws.Range("J24").Value = 279
$ws.Range("K24").Value = 14.695340501792
$ws.Range("L24").Value = 63.265306122449
$ws.Range("M24").Value = 24.513618677042
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 85.185185185185
$ws.Range("I25").Value = 181
$ws.Range("J25").Value = 144
$ws.Range("K25").Value = 25.694444444444
$ws.Range("L25").Value = 57.391304347826
$ws.Range("M25").Value = -33.699633699633
$ws.Range("I26").Value = 16
$ws.Range("K26").Value = 14.285714285714
$ws.Range("L26").Value = 100
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 11
$ws.Range("K27").Value = -27.272727272727
$ws.Range("L27").Value = -60
$ws.Range("G28").Value = 4
$ws.Range("M28").Value = -70
$ws.Range("N28").Value = -92.682926829268
$ws.Range("G29").Value = 3
$ws.Range("M29").Value = -66.666666666666
$ws.Range("N29").Value = -93.150684931506

# --- Cells converting from numeric to text (shared-string) values ---
# Using a leading apostrophe forces text entry, then PasteSpecial(xlPasteFormats)
# from a donor cell that already carries style 14 (t="s") restores the original
# numFmt/style (avoids the auto "quote prefix" style Excel would otherwise add).
$ws.Range("C20").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$ws.Range("A22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("A22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C27").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$excel.CutCopyMode = 0
